# Auto-generated Excel COM-interop script
# Applies scheduled Market Board price-refresh updates to Atomos Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 953.2
$ws.Range("I11").Value = 953.2
$ws.Range("K11").Value = 953.2
$ws.Range("M11").Value = -813.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4513.263
$ws.Range("I80").Value = 393.625
$ws.Range("J80").Value = 7509.364
$ws.Range("K80").Value = 1180.875
$ws.Range("L80").Value = 22528.092
$ws.Range("M80").Value = -182.875
$ws.Range("N80").Value = -24524.092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 4513.263
$ws.Range("I83").Value = 393.625
$ws.Range("J83").Value = 7509.364
$ws.Range("K83").Value = 3542.625
$ws.Range("L83").Value = 67584.276
$ws.Range("M83").Value = 1449.375
$ws.Range("N83").Value = -77568.276

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 587.6667
$ws.Range("I106").Value = 588
$ws.Range("J106").Value = 586
$ws.Range("K106").Value = 588
$ws.Range("L106").Value = 586
$ws.Range("M106").Value = 43
$ws.Range("N106").Value = -1848

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3663.7
$ws.Range("I113").Value = 3281.1667
$ws.Range("J113").Value = 4237.5
$ws.Range("K113").Value = 3281.1667
$ws.Range("L113").Value = 4237.5
$ws.Range("M113").Value = -27.16670000000022
$ws.Range("N113").Value = -10745.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 557380.1
$ws.Range("I141").Value = 1457.6471
$ws.Range("J141").Value = 1907477.6
$ws.Range("K141").Value = 4372.9413
$ws.Range("L141").Value = 5722432.800000001
$ws.Range("M141").Value = 807.0587000000005
$ws.Range("N141").Value = -5732792.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 605.375
$ws.Range("I97").Value = 641.2308
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 641.2308
$ws.Range("L97").Value = 450
$ws.Range("M97").Value = -145.2308
$ws.Range("N97").Value = -1442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 11907184
$ws.Range("I132").Value = 13700800
$ws.Range("J132").Value = 4090.7273
$ws.Range("K132").Value = 41102400
$ws.Range("L132").Value = 12272.1819
$ws.Range("M132").Value = -41099870
$ws.Range("N132").Value = -17332.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3969.0667
$ws.Range("I134").Value = 2821.4546
$ws.Range("J134").Value = 7125
$ws.Range("K134").Value = 8464.363799999999
$ws.Range("L134").Value = 21375
$ws.Range("M134").Value = -5929.363799999999
$ws.Range("N134").Value = -26445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3944.4
$ws.Range("I16").Value = 2222
$ws.Range("J16").Value = 4375
$ws.Range("K16").Value = 2222
$ws.Range("L16").Value = 4375
$ws.Range("M16").Value = -1935
$ws.Range("N16").Value = -4949

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2380.9827
$ws.Range("I31").Value = 1458.5428
$ws.Range("J31").Value = 3784.6956
$ws.Range("K31").Value = 1458.5428
$ws.Range("L31").Value = 3784.6956
$ws.Range("M31").Value = -1163.5428
$ws.Range("N31").Value = -4374.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2380.9827
$ws.Range("I34").Value = 1458.5428
$ws.Range("J34").Value = 3784.6956
$ws.Range("K34").Value = 1458.5428
$ws.Range("L34").Value = 3784.6956
$ws.Range("M34").Value = -1256.5428
$ws.Range("N34").Value = -4188.6956

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19233560
$ws.Range("I58").Value = 1711.6666
$ws.Range("J58").Value = 35718000
$ws.Range("K58").Value = 1711.6666
$ws.Range("L58").Value = 35718000
$ws.Range("M58").Value = -1508.6666
$ws.Range("N58").Value = -35718406

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3944.4
$ws.Range("I113").Value = 2222
$ws.Range("J113").Value = 4375
$ws.Range("K113").Value = 2222
$ws.Range("L113").Value = 4375
$ws.Range("M113").Value = -52
$ws.Range("N113").Value = -8715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4631.125
$ws.Range("I132").Value = 3275
$ws.Range("K132").Value = 9825
$ws.Range("M132").Value = -7295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 19233560
$ws.Range("I136").Value = 1711.6666
$ws.Range("J136").Value = 35718000
$ws.Range("K136").Value = 5134.9998
$ws.Range("L136").Value = 107154000
$ws.Range("M136").Value = -2584.9998
$ws.Range("N136").Value = -107159100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3250
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 9750
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -8814
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3250
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 29250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -24570
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 25002.334
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 25002.334
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 25002.334
$ws.Range("N24").Value = -25348.334
$ws.Range("M24").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 23200
$ws.Range("J114").Value = 23200
$ws.Range("L114").Value = 23200
$ws.Range("N114").Value = -31878

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7069.769
$ws.Range("I122").Value = 9475
$ws.Range("J122").Value = 6000.778
$ws.Range("K122").Value = 28425
$ws.Range("L122").Value = 18002.334
$ws.Range("M122").Value = -25975
$ws.Range("N122").Value = -22902.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2653.25
$ws.Range("I7").Value = 1847.25
$ws.Range("J7").Value = 3056.25
$ws.Range("K7").Value = 1847.25
$ws.Range("L7").Value = 3056.25
$ws.Range("M7").Value = -1735.25
$ws.Range("N7").Value = -3280.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2918.4211
$ws.Range("I122").Value = 2514.2856
$ws.Range("J122").Value = 4050
$ws.Range("K122").Value = 7542.8568
$ws.Range("L122").Value = 12150
$ws.Range("M122").Value = -5092.8568
$ws.Range("N122").Value = -17050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2653.25
$ws.Range("I126").Value = 1847.25
$ws.Range("J126").Value = 3056.25
$ws.Range("K126").Value = 5541.75
$ws.Range("L126").Value = 9168.75
$ws.Range("M126").Value = -3071.75
$ws.Range("N126").Value = -14108.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3506.8572
$ws.Range("I132").Value = 2321.0833
$ws.Range("J132").Value = 5087.8887
$ws.Range("K132").Value = 6963.249899999999
$ws.Range("L132").Value = 15263.6661
$ws.Range("M132").Value = -4433.249899999999
$ws.Range("N132").Value = -20323.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1780.7727
$ws.Range("I136").Value = 1510.3103
$ws.Range("J136").Value = 2303.6667
$ws.Range("K136").Value = 4530.9309
$ws.Range("L136").Value = 6911.000100000001
$ws.Range("M136").Value = -1980.9309
$ws.Range("N136").Value = -12011.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10576
$ws.Range("M15").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 502543.5
$ws.Range("I122").Value = 668264.7
$ws.Range("J122").Value = 5380
$ws.Range("K122").Value = 2004794.1
$ws.Range("L122").Value = 16140
$ws.Range("M122").Value = -2002344.1
$ws.Range("N122").Value = -21040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3227901.5
$ws.Range("I126").Value = 1712.9445
$ws.Range("J126").Value = 7694932
$ws.Range("K126").Value = 5138.833500000001
$ws.Range("L126").Value = 23084796
$ws.Range("M126").Value = -2668.833500000001
$ws.Range("N126").Value = -23089736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16901.572
$ws.Range("I132").Value = 2045.7391
$ws.Range("J132").Value = 45375.25
$ws.Range("K132").Value = 6137.2173
$ws.Range("L132").Value = 136125.75
$ws.Range("M132").Value = -3607.2173
$ws.Range("N132").Value = -141185.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2642.9
$ws.Range("I136").Value = 1291.6364
$ws.Range("J136").Value = 4294.4443
$ws.Range("K136").Value = 3874.9092
$ws.Range("L136").Value = 12883.3329
$ws.Range("M136").Value = -1324.9092
$ws.Range("N136").Value = -17983.3329
